$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.594.86"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "1.630.31"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'212.67"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").Value = "'18.90"
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("E11").Value = "  +3.26%  "
$ws.Range("D12").Value = "1.858.96"
$ws.Range("D13").Value = "1.638.39"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").Value = "'4.06"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").Value = "'0.523"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").Value = "26.598.78"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "'62.91"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'1.00"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'208.74"
$ws.Range("E20").Value = "  +3.94%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'9.42"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").Value = "'6.16"
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("D25").Value = "'146.32"
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("E28").Value = "  +4.50%  "
$ws.Range("D29").Value = "'15.35"
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("E30").Value = "  +3.30%  "
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("D34").Value = "'1.50"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.0171"
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.163.01"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").Value = "'0.807"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").Value = "'0.502"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").Value = "'2.31"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "1.770.67"
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("D45").Value = "'92.39"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").Value = "'54.40"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").Value = "'0.0512"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").Value = "'7.51"
$ws.Range("E50").Value = "  +3.69%  "
$ws.Range("E51").Value = "  +0.06%  "